$wb = $excel.ActiveWorkbook

# --- Rename first worksheet from "Rus" to "Arm" ---
$ws = $wb.Worksheets.Item(1)
$ws.Name = "Arm"

# --- Append a new data row (row 3) with a fresh "Run" record ---
$ws.Range("A3").Value = "12.09.2022, 20:01:58"
$ws.Range("B3").Value = "23-3-2024,14-5-2024,18-5-2024,21-5-2024,25-5-2024"
$ws.Range("C3").Value = "Pass"
$ws.Range("D3").Value = $true
$ws.Range("E3").Value = 2

# --- Keep selection/view pointing at the freshly edited area ---
$ws.Range("B21").Select() | Out-Null
